$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.10851146413066
$ws.Cells.Item(2, 3).Value = 11.05391165979971
$ws.Cells.Item(2, 4).Value = 11.48359070204429
$ws.Cells.Item(2, 6).Value = 29.49622051426013
$ws.Cells.Item(2, 7).Value = 3.626073317721493
$ws.Cells.Item(2, 9).Value = 18.81552094442745
$ws.Cells.Item(2, 10).Value = 11.30401753026426
$ws.Cells.Item(2, 15).Value = 21.062584078276
$ws.Cells.Item(3, 2).Value = 15.34610451830353
$ws.Cells.Item(3, 3).Value = 10.4013205504982
$ws.Cells.Item(3, 4).Value = 11.37995789418844
$ws.Cells.Item(3, 6).Value = 29.56484383084852
$ws.Cells.Item(3, 7).Value = 3.628327604257554
$ws.Cells.Item(3, 9).Value = 18.98843690469477
$ws.Cells.Item(3, 10).Value = 11.28528854685482
$ws.Cells.Item(3, 15).Value = 21.18141356598817
$ws.Cells.Item(4, 2).Value = 14.85790877311013
$ws.Cells.Item(4, 3).Value = 9.978252396925766
$ws.Cells.Item(4, 4).Value = 11.31766441546204
$ws.Cells.Item(4, 6).Value = 29.61747740993712
$ws.Cells.Item(4, 7).Value = 3.629783771727726
$ws.Cells.Item(4, 9).Value = 19.10021739879007
$ws.Cells.Item(4, 10).Value = 11.27611607421313
$ws.Cells.Item(4, 15).Value = 21.26136183060021
$ws.Cells.Item(5, 2).Value = 14.65413578693503
$ws.Cells.Item(5, 3).Value = 9.800303304913664
$ws.Cells.Item(5, 4).Value = 11.29263831841242
$ws.Cells.Item(5, 6).Value = 29.64155228747178
$ws.Cells.Item(5, 7).Value = 3.630395342900458
$ws.Cells.Item(5, 9).Value = 19.14718089753013
$ws.Cells.Item(5, 10).Value = 11.27296578279413
$ws.Cells.Item(5, 15).Value = 21.29568891017075
$ws.Cells.Item(6, 2).Value = 14.62001484984883
$ws.Cells.Item(6, 3).Value = 9.770422830148599
$ws.Cells.Item(6, 4).Value = 11.2885050866974
$ws.Cells.Item(6, 6).Value = 29.64570809723799
$ws.Cells.Item(6, 7).Value = 3.630497993066865
$ws.Cells.Item(6, 9).Value = 19.15506448071311
$ws.Cells.Item(6, 10).Value = 11.27247822702303
$ws.Cells.Item(6, 15).Value = 21.30149416867182
$ws.Cells.Item(7, 2).Value = 14.85517985801721
$ws.Cells.Item(7, 3).Value = 9.975874844722773
$ws.Cells.Item(7, 4).Value = 11.31732542170703
$ws.Cells.Item(7, 6).Value = 29.61779147761073
$ws.Cells.Item(7, 7).Value = 3.629791945937292
$ws.Cells.Item(7, 9).Value = 19.10084504562998
$ws.Cells.Item(7, 10).Value = 11.27607120646368
$ws.Cells.Item(7, 15).Value = 21.2618177143465
$ws.Cells.Item(8, 2).Value = 15.84992073104715
$ws.Cells.Item(8, 3).Value = 10.83358468673863
$ws.Cells.Item(8, 4).Value = 11.44759431083397
$ws.Cells.Item(8, 6).Value = 29.51769501292349
$ws.Cells.Item(8, 7).Value = 3.626835682507378
$ws.Cells.Item(8, 9).Value = 18.87397903306866
$ws.Cells.Item(8, 10).Value = 11.29707832563793
$ws.Cells.Item(8, 15).Value = 21.10210086686414
$ws.Cells.Item(9, 2).Value = 17.63336675814609
$ws.Cells.Item(9, 3).Value = 12.33537166884468
$ws.Cells.Item(9, 4).Value = 11.71264020430491
$ws.Cells.Item(9, 6).Value = 29.40525861739872
$ws.Cells.Item(9, 7).Value = 3.621607235762687
$ws.Cells.Item(9, 9).Value = 18.47353495093181
$ws.Cells.Item(9, 10).Value = 11.35661502875494
$ws.Cells.Item(9, 15).Value = 20.84474460269296
$ws.Cells.Item(10, 2).Value = 18.83305440302279
$ws.Cells.Item(10, 3).Value = 13.32635314936742
$ws.Cells.Item(10, 4).Value = 11.91181522690055
$ws.Cells.Item(10, 6).Value = 29.37443044667602
$ws.Cells.Item(10, 7).Value = 3.618108820619738
$ws.Cells.Item(10, 9).Value = 18.20632556285357
$ws.Cells.Item(10, 10).Value = 11.41134506138749
$ws.Cells.Item(10, 15).Value = 20.69027334146211
$ws.Cells.Item(11, 2).Value = 19.35345684509262
$ws.Cells.Item(11, 3).Value = 13.7524589111671
$ws.Cells.Item(11, 4).Value = 12.0030873806089
$ws.Cells.Item(11, 6).Value = 29.37175836120711
$ws.Cells.Item(11, 7).Value = 3.616590950290865
$ws.Cells.Item(11, 9).Value = 18.09060950823104
$ws.Cells.Item(11, 10).Value = 11.43857848938953
$ws.Cells.Item(11, 15).Value = 20.62763438082174
$ws.Cells.Item(12, 2).Value = 19.54678857899937
$ws.Cells.Item(12, 3).Value = 13.91024628801863
$ws.Cells.Item(12, 4).Value = 12.03772003170617
$ws.Cells.Item(12, 6).Value = 29.37238516673973
$ws.Cells.Item(12, 7).Value = 3.61602669074961
$ws.Cells.Item(12, 9).Value = 18.04762958367549
$ws.Cells.Item(12, 10).Value = 11.44922201170774
$ws.Cells.Item(12, 15).Value = 20.60502166570505
$ws.Cells.Item(13, 2).Value = 19.5053184250036
$ws.Cells.Item(13, 3).Value = 13.87642293770135
$ws.Cells.Item(13, 4).Value = 12.03025857530519
$ws.Cells.Item(13, 6).Value = 29.37217721119961
$ws.Cells.Item(13, 7).Value = 3.616147746909906
$ws.Cells.Item(13, 9).Value = 18.05684875863307
$ws.Cells.Item(13, 10).Value = 11.44691511260754
$ws.Cells.Item(13, 15).Value = 20.60984231309423
$ws.Cells.Item(14, 2).Value = 19.36943763976146
$ws.Cells.Item(14, 3).Value = 13.76551179653922
$ws.Cells.Item(14, 4).Value = 12.00593537768651
$ws.Cells.Item(14, 6).Value = 29.37177705772348
$ws.Cells.Item(14, 7).Value = 3.616544317696788
$ws.Cells.Item(14, 9).Value = 18.08705670608965
$ws.Cells.Item(14, 10).Value = 11.43944753951636
$ws.Cells.Item(14, 15).Value = 20.62575177035462
$ws.Cells.Item(15, 2).Value = 19.28571816624662
$ws.Cells.Item(15, 3).Value = 13.69711021530708
$ws.Cells.Item(15, 4).Value = 11.99104504707948
$ws.Cells.Item(15, 6).Value = 29.37174551188773
$ws.Cells.Item(15, 7).Value = 3.616788597790701
$ws.Cells.Item(15, 9).Value = 18.10566923057682
$ws.Cells.Item(15, 10).Value = 11.43491635632821
$ws.Cells.Item(15, 15).Value = 20.63564127016678
$ws.Cells.Item(16, 2).Value = 18.79852693266317
$ws.Cells.Item(16, 3).Value = 13.29800746447809
$ws.Cells.Item(16, 4).Value = 11.90586156544707
$ws.Cells.Item(16, 6).Value = 29.37483407839268
$ws.Cells.Item(16, 7).Value = 3.61820949283348
$ws.Cells.Item(16, 9).Value = 18.21400533498076
$ws.Cells.Item(16, 10).Value = 11.40961187687111
$ws.Cells.Item(16, 15).Value = 20.69452135353388
$ws.Cells.Item(17, 2).Value = 18.49309145592402
$ws.Cells.Item(17, 3).Value = 13.04683032734951
$ws.Cells.Item(17, 4).Value = 11.85375644634969
$ws.Cells.Item(17, 6).Value = 29.37964130781092
$ws.Cells.Item(17, 7).Value = 3.61909997119598
$ws.Cells.Item(17, 9).Value = 18.28196116810656
$ws.Cells.Item(17, 10).Value = 11.39468328736669
$ws.Cells.Item(17, 15).Value = 20.73260395925935
$ws.Cells.Item(18, 2).Value = 18.31503198188075
$ws.Cells.Item(18, 3).Value = 12.90003594140695
$ws.Cells.Item(18, 4).Value = 11.82385123948679
$ws.Cells.Item(18, 6).Value = 29.38347451656939
$ws.Cells.Item(18, 7).Value = 3.619619079418041
$ws.Cells.Item(18, 9).Value = 18.32159715770853
$ws.Cells.Item(18, 10).Value = 11.38631701409938
$ws.Cells.Item(18, 15).Value = 20.75522572151861
$ws.Cells.Item(19, 2).Value = 18.25433814709236
$ws.Cells.Item(19, 3).Value = 12.84993541021873
$ws.Cells.Item(19, 4).Value = 11.81373768308361
$ws.Cells.Item(19, 6).Value = 29.38495562263381
$ws.Cells.Item(19, 7).Value = 3.619796032226902
$ws.Cells.Item(19, 9).Value = 18.33511163005348
$ws.Cells.Item(19, 10).Value = 11.38352231715849
$ws.Cells.Item(19, 15).Value = 20.76300801280874
$ws.Cells.Item(20, 2).Value = 18.5258527141048
$ws.Cells.Item(20, 3).Value = 13.07380923415823
$ws.Cells.Item(20, 4).Value = 11.85929666097543
$ws.Cells.Item(20, 6).Value = 29.3790189623932
$ws.Cells.Item(20, 7).Value = 3.619004461568496
$ws.Cells.Item(20, 9).Value = 18.27467027622146
$ws.Cells.Item(20, 10).Value = 11.3962497015
$ws.Cells.Item(20, 15).Value = 20.72847565342581
$ws.Cells.Item(21, 2).Value = 19.40945107460378
$ws.Cells.Item(21, 3).Value = 13.79818608024918
$ws.Cells.Item(21, 4).Value = 12.01307800055146
$ws.Cells.Item(21, 6).Value = 29.37185007806428
$ws.Cells.Item(21, 7).Value = 3.61642755003611
$ws.Cells.Item(21, 9).Value = 18.0781611253435
$ws.Cells.Item(21, 10).Value = 11.44163201189327
$ws.Cells.Item(21, 15).Value = 20.62104864620713
$ws.Cells.Item(22, 2).Value = 19.9651421959559
$ws.Cells.Item(22, 3).Value = 14.25079743477249
$ws.Cells.Item(22, 4).Value = 12.11397830725583
$ws.Cells.Item(22, 6).Value = 29.37671787867423
$ws.Cells.Item(22, 7).Value = 3.614804711857669
$ws.Cells.Item(22, 9).Value = 17.95462267323239
$ws.Cells.Item(22, 10).Value = 11.47321751131506
$ws.Cells.Item(22, 15).Value = 20.55729762658141
$ws.Cells.Item(23, 2).Value = 19.67057874275707
$ws.Cells.Item(23, 3).Value = 14.0111388848265
$ws.Cells.Item(23, 4).Value = 12.06009828310516
$ws.Cells.Item(23, 6).Value = 29.37324410819021
$ws.Cells.Item(23, 7).Value = 3.615665258257192
$ws.Cells.Item(23, 9).Value = 18.02010996888198
$ws.Cells.Item(23, 10).Value = 11.45618536827835
$ws.Cells.Item(23, 15).Value = 20.59072849899053
$ws.Cells.Item(24, 2).Value = 18.51104900641892
$ws.Cells.Item(24, 3).Value = 13.06161951116703
$ws.Cells.Item(24, 4).Value = 11.85679176998411
$ws.Cells.Item(24, 6).Value = 29.37929699386069
$ws.Cells.Item(24, 7).Value = 3.619047619162094
$ws.Cells.Item(24, 9).Value = 18.27796472129964
$ws.Cells.Item(24, 10).Value = 11.39554085131171
$ws.Cells.Item(24, 15).Value = 20.73033979486093
$ws.Cells.Item(25, 2).Value = 17.1697937858121
$ws.Cells.Item(25, 3).Value = 11.94869647948781
$ws.Cells.Item(25, 4).Value = 11.64005495285156
$ws.Cells.Item(25, 6).Value = 29.42662452301256
$ws.Cells.Item(25, 7).Value = 3.622961173594798
$ws.Cells.Item(25, 9).Value = 18.57711842782625
$ws.Cells.Item(25, 10).Value = 11.33856324297035
$ws.Cells.Item(25, 15).Value = 20.90832911436895
